$d = $word.ActiveDocument

function Get-ParaRangeAfter($startRange, $searchText) {
    $scan = $d.Range($startRange.End, $d.Content.End)
    $ok = $scan.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "Not found (after): $searchText" }
    return $scan.Paragraphs(1).Range
}

function Get-ParaRange($searchText) {
    $scan = $d.Content
    $ok = $scan.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "Not found: $searchText" }
    return $scan.Paragraphs(1).Range
}

# Remove bookmarkStart (_GoBack) from the title paragraph
$p1 = Get-ParaRange "SETTING UP YOUR BREADBOARD AND CONNECTING DEVICES TO IT"
$p1Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="38498547" w14:textId="77777777" w:rsidR="00B6030F" w:rsidRPr="00240701" w:rsidRDefault="00240701"><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr><w:r w:rsidRPr="00240701"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">LAB </w:t></w:r><w:r w:rsidR="00F36588"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>1-1</w:t></w:r><w:r w:rsidRPr="00240701"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="007476FF"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>SETTING UP YOUR BREADBOARD AND CONNECTING DEVICES TO IT</w:t></w:r></w:p>
"@
$p1.InsertXML($p1Xml)

# Delete the "Adel was here " run entirely
$p2 = Get-ParaRange "Adel was here"
$p2Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="778ABA53" w14:textId="6FF4F38E" w:rsidR="00240701" w:rsidRPr="00240701" w:rsidRDefault="00A217F5"><w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>
"@
$p2.InsertXML($p2Xml)

# Split "Adafruit T " run, wrapping "Adafruit" in proofErr spellStart/spellEnd (first occurrence)
$p3 = Get-ParaRange "In this lab, you will setup"
$p3Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="44B24BC6" w14:textId="77777777" w:rsidR="00B90E41" w:rsidRPr="00240701" w:rsidRDefault="00B90E41" w:rsidP="007476FF"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00240701"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">In this lab, you will setup </w:t></w:r><w:r w:rsidR="007476FF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>your breadboard to be connected to</w:t></w:r><w:r w:rsidRPr="00240701"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the Raspberry Pi using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="007476FF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Adafruit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="007476FF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> T </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="007476FF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>CobblerPlus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="007476FF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Kit and connect a</w:t></w:r><w:r w:rsidR="00155037"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>n</w:t></w:r><w:r w:rsidR="007476FF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> LED light and Push Button switch to it and write Python programs to make these devices work.</w:t></w:r></w:p>
"@
$p3.InsertXML($p3Xml)

# Split "Adafruit T Cobbler Plus and Breakout Cable" run, wrapping "Adafruit" in proofErr (second occurrence)
$p4 = Get-ParaRange "Adafruit T Cobbler Plus and Breakout Cable"
$p4Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4133CE9E" w14:textId="77777777" w:rsidR="00B90E41" w:rsidRDefault="00155037" w:rsidP="00B90E41"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Adafruit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> T Cobbler Plus and Breakout Cable</w:t></w:r></w:p>
"@
$p4.InsertXML($p4Xml)

# Split "Turn on your Pi and wait a few mins" run, wrapping "mins" in proofErr
$p5 = Get-ParaRange "Turn on your Pi and wait a few mins"
$p5Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="26EFE7BC" w14:textId="77777777" w:rsidR="001922E4" w:rsidRDefault="001922E4" w:rsidP="001922E4"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="9"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Turn on your Pi and wait a few </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>mins</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@
$p5.InsertXML($p5Xml)

# Split "def blinking(pin):" run, wrapping "def" in proofErr
$p6 = Get-ParaRange "def blinking(pin):"
$p6Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1CD5C872" w14:textId="77777777" w:rsidR="00CB0C88" w:rsidRPr="00CB0C88" w:rsidRDefault="00CB0C88" w:rsidP="00CB0C88"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00CB0C88"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>def</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00CB0C88"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> blinking(pin):</w:t></w:r></w:p>
"@
$p6.InsertXML($p6Xml)

# Add bookmarkStart/bookmarkEnd (_GoBack) after the "PUSH_BUTTON = 13" run
$p7 = Get-ParaRange "PUSH_BUTTON = 13"
$p7Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="10D38826" w14:textId="77777777" w:rsidR="002445D3" w:rsidRPr="002445D3" w:rsidRDefault="00C62979" w:rsidP="002445D3"><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>PUSH_BUTTON = 13</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@
$p7.InsertXML($p7Xml)

# Remove the old bookmarkEnd (_GoBack) from the final "GPIO.cleanup()" paragraph
$anchor = Get-ParaRange "PUSH_BUTTON = 13"
$p8 = Get-ParaRangeAfter $anchor "GPIO.cleanup()"
$p8Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4B19D7B7" w14:textId="77777777" w:rsidR="00C37860" w:rsidRPr="002445D3" w:rsidRDefault="002445D3" w:rsidP="002445D3"><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:sz w:val="32"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="002445D3"><w:rPr><w:sz w:val="28"/></w:rPr><w:t>GPIO.cleanup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="002445D3"><w:rPr><w:sz w:val="28"/></w:rPr><w:t>()</w:t></w:r></w:p>
"@
$p8.InsertXML($p8Xml)

Write-Host "All edits applied"